$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B41: convert from text "3" to numeric 3
$ws.Range("B41").Value = 3

# Add new row 42 (mirrors old row 41's content pattern, with B42 kept as text "3")
$ws.Range("A42").Value = "Ruilin"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "3"
$ws.Range("B42").Style = "Normal"

$ws.Range("C42").Value = "无"
$ws.Range("D42").Value = "CRT"
$ws.Range("E42").Value = "MET"
$ws.Range("F42").Value = "295c014b-37cb-453e-93b8-ae293d0d968b"
$ws.Range("G42").Value = "BkiIkBJ0b_annotated.xlsx"
$ws.Range("H42").Value = "The other part of the criticism that we use a ""straw man"" is again wrong because we do not intend to show pathology with Mirowski et al. paper, experiments or claims."
